# Add loop to check repeat user or new user
# -> append two new user rows to the "Users" sheet (rows 6 and 7),
#    extending the data below the existing Katja/Tanja/Julja/Tanja rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A6").Value = "Yulja"
$ws.Range("B6").Value = "yulja@gmail.com"
$ws.Range("C6").Value = -7
$ws.Range("C6").Style = "Normal"

$ws.Range("A7").Value = "tanja"
$ws.Range("B7").Value = "tanja@gmail.com"
$ws.Range("C7").Value = -9
$ws.Range("C7").Style = "Normal"
